$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.402.58'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").Value = '1.840.82'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("E4").Value = '  +1.37%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '315.11'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("E6").Value = '  +1.17%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4765'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.59%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3706'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.32%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07471'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.32%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.8868'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.73%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '20.51'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").Value = '1.875.72'
$ws.Range("E12").Value = '  +3.03%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.07360'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.13%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.490'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.99%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '93.27'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.49%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '6.596'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.14%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.015'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.27%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008856'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("E19").Value = '  +1.25%  '
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D21").Value = '27.421.68'
$ws.Range("E21").Value = '  +1.75%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.356'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.59%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.75'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("D24").Value = '2.075.78'
$ws.Range("E24").Value = '  +1.30%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.900'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.41%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '152.68'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.10%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.65'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.67%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.172'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.00%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.278'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '118.22'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.71%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.08996'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.40%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.7610'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.95%  '
$ws.Range("E33").Value = '  +1.22%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.565'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.18%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.948'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("E36").Value = '  +1.29%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.108'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.84%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.05385'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.74%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01967'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.14%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.999'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.69%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '7.315'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.56%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.5367'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.84%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.398'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.79%  '
$ws.Range("E44").Value = '  +0.62%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.571'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.39%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4995'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.38%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.59'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.07%  '
$ws.Range("E48").Value = '  +1.33%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '105.25'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.44%  '
$ws.Range("E50").Value = '  +0.67%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.06328'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.39%  '
